$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells being updated so numeric-looking
# strings (e.g. "216.45", "0.0520") are preserved exactly as text, not coerced
# to floating point numbers by the usual Excel type inference on Range.Value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = "26.784.37"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.648.52"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "216.45"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.878.04"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "1.666.24"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "65.32"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.800.45"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "218.31"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "2.54"
$ws.Range("E21").Value = "  +14.34%  "
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "6.30"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "146.28"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.55"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.280.42"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "0.542"
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").Value = "1.789.41"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "92.11"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "59.95"
$ws.Range("E46").Value = "  +9.04%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("D51").Value = "0.0980"
$ws.Range("E51").Value = "  +2.35%  "
